$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = -0.5635264962161868
$ws.Cells.Item(2, 5).Value = -32.28768987698245
$ws.Cells.Item(3, 1).Value = 0.02638492562814984
$ws.Cells.Item(3, 2).Value = 0.116217378282943
$ws.Cells.Item(3, 3).Value = 2.426565366991384
$ws.Cells.Item(3, 4).Value = -0.5635264962161868
$ws.Cells.Item(3, 5).Value = -32.28768987698245
$ws.Cells.Item(4, 1).Value = 0.05276985125629968
$ws.Cells.Item(4, 2).Value = 0.2324347565658859
$ws.Cells.Item(4, 3).Value = 2.353130733982769
$ws.Cells.Item(4, 4).Value = -0.5621414372011267
$ws.Cells.Item(4, 5).Value = -32.20833184104297
$ws.Cells.Item(5, 1).Value = 0.07915477688444952
$ws.Cells.Item(5, 2).Value = 0.3487537346417531
$ws.Cells.Item(5, 3).Value = 2.279857139288249
$ws.Cells.Item(5, 4).Value = -0.5593691641439307
$ws.Cells.Item(5, 5).Value = -32.04949229520781
$ws.Cells.Item(6, 1).Value = 0.1055397025125994
$ws.Cells.Item(6, 2).Value = 0.4652753998854529
$ws.Cells.Item(6, 3).Value = 2.206906293721485
$ws.Cells.Item(6, 4).Value = -0.5552053535574338
$ws.Cells.Item(6, 5).Value = -31.81092352190965
$ws.Cells.Item(7, 1).Value = 0.1319246281407492
$ws.Cells.Item(7, 2).Value = 0.5820998076693485
$ws.Cells.Item(7, 3).Value = 2.134441253281196
$ws.Cells.Item(7, 4).Value = -0.5496434864500335
$ws.Cells.Item(7, 5).Value = -31.49225201044297
$ws.Cells.Item(8, 1).Value = 0.158309553768899
$ws.Cells.Item(8, 2).Value = 0.6993254473603198
$ws.Cells.Item(8, 3).Value = 2.062627092149053
$ws.Cells.Item(8, 4).Value = -0.5426748073092807
$ws.Cells.Item(8, 5).Value = -31.09297610689698
$ws.Cells.Item(9, 1).Value = 0.1846944793970489
$ws.Cells.Item(9, 2).Value = 0.81704868647458
$ws.Cells.Item(9, 3).Value = 1.991631576005304
$ws.Cells.Item(9, 4).Value = -0.5342882685620018
$ws.Cells.Item(9, 5).Value = -30.61246283195497
$ws.Cells.Item(10, 1).Value = 0.2110794050251987
$ws.Cells.Item(10, 2).Value = 0.9353631853055511
$ws.Cells.Item(10, 3).Value = 1.921625835480351
$ws.Cells.Item(10, 4).Value = -0.5244704597560446
$ws.Cells.Item(10, 5).Value = -30.04994382330725
$ws.Cells.Item(11, 1).Value = 0.2374643306533486
$ws.Cells.Item(11, 2).Value = 1.054359273986733
$ws.Cells.Item(11, 3).Value = 1.852785039300681
$ws.Cells.Item(11, 4).Value = -0.5132055204874867
$ws.Cells.Item(11, 5).Value = -29.40451034674769
$ws.Cells.Item(12, 1).Value = 0.2638492562814984
$ws.Cells.Item(12, 2).Value = 1.174123283500811
$ws.Cells.Item(12, 3).Value = 1.78528906634538
$ws.Cells.Item(12, 4).Value = -0.5004750358525043
$ws.Cells.Item(12, 5).Value = -28.67510730600706
$ws.Cells.Item(13, 1).Value = 0.2902341819096482
$ws.Cells.Item(13, 2).Value = 1.294736821592092
$ws.Cells.Item(13, 3).Value = 1.719323175392752
$ws.Cells.Item(13, 4).Value = -0.4862579129285956
$ws.Cells.Item(13, 5).Value = -27.86052616564839
$ws.Cells.Item(14, 1).Value = 0.3166191075377981
$ws.Cells.Item(14, 2).Value = 1.416275983874329
$ws.Cells.Item(14, 3).Value = 1.655078670788834
$ws.Cells.Item(14, 4).Value = -0.4705302364778881
$ws.Cells.Item(14, 5).Value = -26.95939668347556
$ws.Cells.Item(15, 1).Value = 0.343004033165948
$ws.Cells.Item(15, 2).Value = 1.538810489641165
$ws.Cells.Item(15, 3).Value = 1.592753561590792
$ws.Cells.Item(15, 4).Value = -0.4532651017067332
$ws.Cells.Item(15, 5).Value = -25.97017732836382
$ws.Cells.Item(16, 1).Value = 0.3693889587940978
$ws.Cells.Item(16, 2).Value = 1.662402730971631
$ws.Cells.Item(16, 3).Value = 1.532553210902588
$ws.Cells.Item(16, 4).Value = -0.4344324214997735
$ws.Cells.Item(16, 5).Value = -24.89114423558546
$ws.Cells.Item(17, 1).Value = 0.3957738844222476
$ws.Cells.Item(17, 2).Value = 1.78710672266658
$ws.Cells.Item(17, 3).Value = 1.474690971096157
$ws.Cells.Item(17, 4).Value = -0.4139987050597896
$ws.Cells.Item(17, 5).Value = -23.7203785238073
$ws.Cells.Item(18, 1).Value = 0.4221588100503975
$ws.Cells.Item(18, 2).Value = 1.912966939340635
$ws.Cells.Item(18, 3).Value = 1.419388799359135
$ws.Cells.Item(18, 4).Value = -0.3919268043104899
$ws.Cells.Item(18, 5).Value = -22.45575176504079
$ws.Cells.Item(19, 1).Value = 0.4485437356785473
$ws.Cells.Item(19, 2).Value = 2.040017024613645
$ws.Cells.Item(19, 3).Value = 1.366877846480788
$ws.Cells.Item(19, 4).Value = -0.3681756237376402
$ws.Cells.Item(19, 5).Value = -21.09490935976339
$ws.Cells.Item(20, 1).Value = 0.4749286613066971
$ws.Cells.Item(20, 2).Value = 2.16827835578049
$ws.Cells.Item(20, 3).Value = 1.317399009920297
$ws.Cells.Item(20, 4).Value = -0.3426997885291662
$ws.Cells.Item(20, 5).Value = -19.63525152274704
$ws.Cells.Item(21, 1).Value = 0.501313586934847
$ws.Cells.Item(21, 2).Value = 2.297758445572029
$ws.Cells.Item(21, 3).Value = 1.271203439920067
$ws.Cells.Item(21, 4).Value = -0.3154492648951779
$ws.Cells.Item(21, 5).Value = -18.07391152899801
$ws.Cells.Item(22, 1).Value = 0.5276985125629968
$ws.Cells.Item(22, 2).Value = 2.428449160637006
$ws.Cells.Item(22, 3).Value = 1.228552984636728
$ws.Cells.Item(22, 4).Value = -0.2863689252639304
$ws.Cells.Item(22, 5).Value = -16.4077308013205
$ws.Cells.Item(23, 1).Value = 0.5540834381911467
$ws.Cells.Item(23, 2).Value = 2.560324734159684
$ws.Cells.Item(23, 3).Value = 1.189720556845822
$ws.Cells.Item(23, 4).Value = -0.2553980496081627
$ws.Cells.Item(23, 5).Value = -14.63323033842055
$ws.Cells.Item(24, 1).Value = 0.5804683638192965
$ws.Cells.Item(24, 2).Value = 2.693339547568394
$ws.Cells.Item(24, 3).Value = 1.154990400584764
$ws.Cells.Item(24, 4).Value = -0.2224697523922488
$ws.Cells.Item(24, 5).Value = -12.74657788139631
$ws.Cells.Item(25, 1).Value = 0.6068532894474463
$ws.Cells.Item(25, 2).Value = 2.827425653579025
$ws.Cells.Item(25, 3).Value = 1.124658230946286
$ws.Cells.Item(25, 4).Value = -0.1875103224592058
$ws.Cells.Item(25, 5).Value = -10.74355009204963
$ws.Cells.Item(26, 1).Value = 0.6332382150755962
$ws.Cells.Item(26, 2).Value = 2.962490009856605
$ws.Cells.Item(26, 3).Value = 1.099031213886263
$ws.Cells.Item(26, 4).Value = -0.1504384604875428
$ws.Cells.Item(26, 5).Value = -8.6194888623818
$ws.Cells.Item(27, 1).Value = 0.6596231407037461
$ws.Cells.Item(27, 2).Value = 3.098411389384343
$ws.Cells.Item(27, 3).Value = 1.078427745067921
$ws.Cells.Item(27, 4).Value = -0.1111643952968098
$ws.Cells.Item(27, 5).Value = -6.36925068263114
$ws.Cells.Item(28, 1).Value = 0.6860080663318959
$ws.Cells.Item(28, 2).Value = 3.235036930243101
$ws.Cells.Item(28, 3).Value = 1.063176977050422
$ws.Cells.Item(28, 4).Value = -0.06958885607717469
$ws.Cells.Item(28, 5).Value = -3.98714775436542
$ws.Cells.Item(29, 1).Value = 0.7123929919600457
$ws.Cells.Item(29, 2).Value = 3.372178284001732
$ws.Cells.Item(29, 3).Value = 1.053618032060397
$ws.Cells.Item(29, 4).Value = -0.02560187230936282
$ws.Cells.Item(29, 5).Value = -1.46687923095934
$ws.Cells.Item(30, 1).Value = 0.7387779175881956
$ws.Cells.Item(30, 2).Value = 3.509607318431962
$ws.Cells.Item(30, 3).Value = 1.050098822540203
$ws.Cells.Item(30, 4).Value = 0.02091863361235319
$ws.Cells.Item(30, 5).Value = 1.198549419168341
$ws.Cells.Item(31, 1).Value = 0.7651628432163454
$ws.Cells.Item(31, 2).Value = 3.647051327006321
$ws.Cells.Item(31, 3).Value = 1.052974382848091
$ws.Cells.Item(31, 4).Value = 0.0701085046985388
$ws.Cells.Item(31, 5).Value = 4.016921427199375
$ws.Cells.Item(32, 1).Value = 0.7915477688444952
$ws.Cells.Item(32, 2).Value = 3.784187694956132
$ws.Cells.Item(32, 3).Value = 1.06260459185969
$ws.Cells.Item(32, 4).Value = 0.1221193134558592
$ws.Cells.Item(32, 5).Value = 6.996921258055896
$ws.Cells.Item(33, 1).Value = 0.817932694472645
$ws.Cells.Item(33, 2).Value = 3.920637970091656
$ws.Cells.Item(33, 3).Value = 1.079351136443969
$ws.Cells.Item(33, 4).Value = 0.1771207011639635
$ws.Cells.Item(33, 5).Value = 10.148268641093
$ws.Cells.Item(34, 1).Value = 0.8443176201007949
$ws.Cells.Item(34, 2).Value = 4.055961286937896
$ws.Cells.Item(34, 3).Value = 1.103573528105533
$ws.Cells.Item(34, 4).Value = 0.2353031809311819
$ws.Cells.Item(34, 5).Value = 13.48187917335991
$ws.Cells.Item(35, 1).Value = 0.8707025457289448
$ws.Cells.Item(35, 2).Value = 4.189647096271455
$ws.Cells.Item(35, 3).Value = 1.135623937196612
$ws.Cells.Item(35, 4).Value = 0.2968815186016825
$ws.Cells.Item(35, 5).Value = 17.01005803131105
$ws.Cells.Item(36, 1).Value = 0.8970874713570947
$ws.Cells.Item(36, 2).Value = 4.321107160769394
$ws.Cells.Item(36, 3).Value = 1.175840547972166
$ws.Cells.Item(36, 4).Value = 0.362098839792845
$ws.Cells.Item(36, 5).Value = 20.74673528671377
$ws.Cells.Item(37, 1).Value = 0.9234723969852444
$ws.Cells.Item(37, 2).Value = 4.449666794120136
$ws.Cells.Item(37, 3).Value = 1.224539059374445
$ws.Cells.Item(37, 4).Value = 0.4312316577879921
$ws.Cells.Item(37, 5).Value = 24.70775398368176
$ws.Cells.Item(38, 1).Value = 0.9498573226133943
$ws.Cells.Item(38, 2).Value = 4.57455535006764
$ws.Cells.Item(38, 3).Value = 1.282001855455242
$ws.Cells.Item(38, 4).Value = 0.504596080855312
$ws.Cells.Item(38, 5).Value = 28.91122579185142
$ws.Cells.Item(39, 1).Value = 0.9762422482415442
$ws.Cells.Item(39, 2).Value = 4.694896016334018
$ws.Cells.Item(39, 3).Value = 1.348464238670421
$ws.Cells.Item(39, 4).Value = 0.5825555465024905
$ws.Cells.Item(39, 5).Value = 33.37797414652987
$ws.Cells.Item(40, 1).Value = 1.002627173869694
$ws.Cells.Item(40, 2).Value = 4.809695046827148
$ws.Cells.Item(40, 3).Value = 1.424096949437092
$ws.Cells.Item(40, 4).Value = 0.6655305558771523
$ws.Cells.Item(40, 5).Value = 38.13209198875643
$ws.Cells.Item(41, 1).Value = 1.029012099497844
$ws.Cells.Item(41, 2).Value = 4.917830690585833
$ws.Cells.Item(41, 3).Value = 1.508983973704486
$ws.Cells.Item(41, 4).Value = 0.7540110620123706
$ws.Cells.Item(41, 5).Value = 43.20165155948582
$ws.Cells.Item(42, 1).Value = 1.055397025125994
$ws.Cells.Item(42, 2).Value = 5.018042273708278
$ws.Cells.Item(42, 3).Value = 1.603094350115236
$ws.Cells.Item(42, 4).Value = 0.8485724293084896
$ws.Cells.Item(42, 5).Value = 48.61961881053986
